$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 3.883322333333334
$ws.Range("H2").Value = 11.649967
$ws.Range("I2").Value = 0.2039370428985499
$ws.Range("J2").Value = 0.2039370428985498
$ws.Range("M2").Value = 57.65261933333333
$ws.Range("N2").Value = 172.957858
$ws.Range("O2").Value = 0.6817060950001529
$ws.Range("P2").Value = 0.6817060950001529
$ws.Range("Q2").Value = 223.8837042322984
$ws.Range("R2").Value = 2014.953338090686
$ws.Range("S2").Value = 0.1390251251402491
$ws.Range("T2").Value = 0.1390251251402491
$ws.Range("G3").Value = 3.883322333333334
$ws.Range("H3").Value = 11.649967
$ws.Range("I3").Value = 0.2039370428985499
$ws.Range("J3").Value = 0.2039370428985498
$ws.Range("O3").Value = 0.1019529789289588
$ws.Range("P3").Value = 0.1019529789289588
$ws.Range("Q3").Value = 33.48306659943778
$ws.Range("R3").Value = 301.3475993949399
$ws.Range("S3").Value = 0.02079198903747003
$ws.Range("T3").Value = 0.02079198903747002
$ws.Range("G4").Value = 3.883322333333334
$ws.Range("H4").Value = 11.649967
$ws.Range("I4").Value = 0.2039370428985499
$ws.Range("J4").Value = 0.2039370428985498
$ws.Range("M4").Value = 2.790736
$ws.Range("N4").Value = 8.372208000000001
$ws.Range("O4").Value = 0.0329987043561157
$ws.Range("P4").Value = 0.0329987043561157
$ws.Range("Q4").Value = 10.83732743523733
$ws.Range("R4").Value = 97.535946917136
$ws.Range("S4").Value = 0.006729658185869732
$ws.Range("T4").Value = 0.00672965818586973
$ws.Range("G5").Value = 3.883322333333334
$ws.Range("H5").Value = 11.649967
$ws.Range("I5").Value = 0.2039370428985499
$ws.Range("J5").Value = 0.2039370428985498
$ws.Range("M5").Value = 15.50544933333333
$ws.Range("N5").Value = 46.516348
$ws.Range("O5").Value = 0.1833422217147727
$ws.Range("P5").Value = 0.1833422217147727
$ws.Range("Q5").Value = 60.21265768450178
$ws.Range("R5").Value = 541.913919160516
$ws.Range("S5").Value = 0.03739027053496104
$ws.Range("T5").Value = 0.03739027053496103
$ws.Range("G6").Value = 9.654910333333332
$ws.Range("I6").Value = 0.5070384824688307
$ws.Range("J6").Value = 0.5070384824688307
$ws.Range("M6").Value = 57.65261933333333
$ws.Range("N6").Value = 172.957858
$ws.Range("O6").Value = 0.6817060950001529
$ws.Range("P6").Value = 0.6817060950001529
$ws.Range("Q6").Value = 556.6308701451329
$ws.Range("R6").Value = 5009.677831306197
$ws.Range("S6").Value = 0.34565122389863
$ws.Range("T6").Value = 0.34565122389863
$ws.Range("G7").Value = 9.654910333333332
$ws.Range("I7").Value = 0.5070384824688307
$ws.Range("J7").Value = 0.5070384824688307
$ws.Range("O7").Value = 0.1019529789289588
$ws.Range("P7").Value = 0.1019529789289588
$ws.Range("Q7").Value = 83.24727590282443
$ws.Range("R7").Value = 749.2254831254198
$ws.Range("S7").Value = 0.05169408371931595
$ws.Range("T7").Value = 0.05169408371931595
$ws.Range("G8").Value = 9.654910333333332
$ws.Range("I8").Value = 0.5070384824688307
$ws.Range("J8").Value = 0.5070384824688307
$ws.Range("M8").Value = 2.790736
$ws.Range("N8").Value = 8.372208000000001
$ws.Range("O8").Value = 0.0329987043561157
$ws.Range("P8").Value = 0.0329987043561157
$ws.Range("Q8").Value = 26.94430584400533
$ws.Range("R8").Value = 242.498752596048
$ws.Range("S8").Value = 0.0167316129801625
$ws.Range("T8").Value = 0.01673161298016249
$ws.Range("G9").Value = 9.654910333333332
$ws.Range("I9").Value = 0.5070384824688307
$ws.Range("J9").Value = 0.5070384824688307
$ws.Range("M9").Value = 15.50544933333333
$ws.Range("N9").Value = 46.516348
$ws.Range("O9").Value = 0.1833422217147727
$ws.Range("P9").Value = 0.1833422217147727
$ws.Range("Q9").Value = 149.7037229913764
$ws.Range("R9").Value = 1347.333506922388
$ws.Range("S9").Value = 0.09296156187072223
$ws.Range("T9").Value = 0.09296156187072223
$ws.Range("G10").Value = 4.652793333333332
$ws.Range("H10").Value = 13.95838
$ws.Range("I10").Value = 0.2443466784802274
$ws.Range("J10").Value = 0.2443466784802274
$ws.Range("M10").Value = 57.65261933333333
$ws.Range("N10").Value = 172.957858
$ws.Range("O10").Value = 0.6817060950001529
$ws.Range("P10").Value = 0.6817060950001529
$ws.Range("Q10").Value = 268.2457228833377
$ws.Range("R10").Value = 2414.21150595004
$ws.Range("S10").Value = 0.1665726200130137
$ws.Range("T10").Value = 0.1665726200130137
$ws.Range("G11").Value = 4.652793333333332
$ws.Range("H11").Value = 13.95838
$ws.Range("I11").Value = 0.2443466784802274
$ws.Range("J11").Value = 0.2443466784802274
$ws.Range("O11").Value = 0.1019529789289588
$ws.Range("P11").Value = 0.1019529789289588
$ws.Range("Q11").Value = 40.1176558835111
$ws.Range("R11").Value = 361.0589029515999
$ws.Range("S11").Value = 0.0249118717624557
$ws.Range("T11").Value = 0.0249118717624557
$ws.Range("G12").Value = 4.652793333333332
$ws.Range("H12").Value = 13.95838
$ws.Range("I12").Value = 0.2443466784802274
$ws.Range("J12").Value = 0.2443466784802274
$ws.Range("M12").Value = 2.790736
$ws.Range("N12").Value = 8.372208000000001
$ws.Range("O12").Value = 0.0329987043561157
$ws.Range("P12").Value = 0.0329987043561157
$ws.Range("Q12").Value = 12.98471785589333
$ws.Range("R12").Value = 116.86246070304
$ws.Range("S12").Value = 0.008063123803567884
$ws.Range("T12").Value = 0.008063123803567883
$ws.Range("G13").Value = 4.652793333333332
$ws.Range("H13").Value = 13.95838
$ws.Range("I13").Value = 0.2443466784802274
$ws.Range("J13").Value = 0.2443466784802274
$ws.Range("M13").Value = 15.50544933333333
$ws.Range("N13").Value = 46.516348
$ws.Range("O13").Value = 0.1833422217147727
$ws.Range("P13").Value = 0.1833422217147727
$ws.Range("Q13").Value = 72.14365128847109
$ws.Range("R13").Value = 649.29286159624
$ws.Range("S13").Value = 0.04479906290119013
$ws.Range("T13").Value = 0.04479906290119013
$ws.Range("G14").Value = 0.8507443333333334
$ws.Range("H14").Value = 2.552233
$ws.Range("I14").Value = 0.04467779615239207
$ws.Range("J14").Value = 0.04467779615239207
$ws.Range("M14").Value = 57.65261933333333
$ws.Range("N14").Value = 172.957858
$ws.Range("O14").Value = 0.6817060950001529
$ws.Range("P14").Value = 0.6817060950001529
$ws.Range("Q14").Value = 49.04763919965711
$ws.Range("R14").Value = 441.428752796914
$ws.Range("S14").Value = 0.03045712594826006
$ws.Range("T14").Value = 0.03045712594826006
$ws.Range("G15").Value = 0.8507443333333334
$ws.Range("H15").Value = 2.552233
$ws.Range("I15").Value = 0.04467779615239207
$ws.Range("J15").Value = 0.04467779615239207
$ws.Range("O15").Value = 0.1019529789289588
$ws.Range("P15").Value = 0.1019529789289588
$ws.Range("Q15").Value = 7.335350178784444
$ws.Range("R15").Value = 66.01815160906
$ws.Range("S15").Value = 0.004555034409717147
$ws.Range("T15").Value = 0.004555034409717147
$ws.Range("G16").Value = 0.8507443333333334
$ws.Range("H16").Value = 2.552233
$ws.Range("I16").Value = 0.04467779615239207
$ws.Range("J16").Value = 0.04467779615239207
$ws.Range("M16").Value = 2.790736
$ws.Range("N16").Value = 8.372208000000001
$ws.Range("O16").Value = 0.0329987043561157
$ws.Range("P16").Value = 0.0329987043561157
$ws.Range("Q16").Value = 2.374202837829334
$ws.Range("R16").Value = 21.367825540464
$ws.Range("S16").Value = 0.00147430938651559
$ws.Range("T16").Value = 0.001474309386515589
$ws.Range("G17").Value = 0.8507443333333334
$ws.Range("H17").Value = 2.552233
$ws.Range("I17").Value = 0.04467779615239207
$ws.Range("J17").Value = 0.04467779615239207
$ws.Range("O17").Value = 0.1019529789289588
$ws.Range("P17").Value = 0.1019529789289588
$ws.Range("Q17").Value = 13.19117315612045
$ws.Range("R17").Value = 118.720558405084
$ws.Range("S17").Value = 0.008191326407899284
$ws.Range("T17").Value = 0.008191326407899284
